$d = $word.ActiveDocument

# 1. Remove the old _GoBack bookmark (currently sitting right before "Error-Page")
$d.Bookmarks("_GoBack").Delete()

# 2. Split the "Wir haben uns entschieden..." paragraph: turn the manual line break before
#    "Die bank haben wir..." into a real paragraph break, rewrite the start of that sentence
#    ("Die bank" -> "Das Bank-Objekt selbst"), and fix "diese ... dieselbe" -> "dieses ... dasselbe".
$findText = "abgefragt wird.^lDie bank haben wir in einem Singleton gespeichert, damit diese immer dieselbe bleibt."
$replText = "abgefragt wird.^pDas Bank-Objekt selbst haben wir in einem Singleton gespeichert, damit dieses immer dasselbe bleibt."
$d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replText, 2) | Out-Null

# 3. Re-insert the _GoBack bookmark in its new spot: right before "haben wir in einem Singleton"
#    (i.e. immediately after "Das Bank-Objekt selbst ").
$full = $d.Content.Text
$pos = $full.IndexOf("haben wir in einem Singleton")
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 4. Insert the new paragraph about the HtmlPrinter class right before the
#    "Eine kleine dynamische Error-Page..." paragraph.
$full = $d.Content.Text
$pos2 = $full.IndexOf("Eine kleine dynamische")
$r = $d.Range($pos2, $pos2)
$r.InsertBefore("Um die Aufgaben etwas zu trennen und die Servlet-Klasse nicht zu überfüllen, übernimmt die Klasse HtmlPrinter die Ausgabe, die bei einem doGet vorgenommen werden muss.`r")

# 5. Normalize the "Eine kleine dynamische " / "Error-Page..." runs (bookmark that used to split
#    them is gone now) back into a single run.
$d.Content.Find.Execute("Eine kleine dynamische Error-Page", $true, $false, $false, $false, $false, $true, 1, $false, "Eine kleine dynamische Error-Page", 2) | Out-Null
